# Generate Report for Archive
#
# 1) Update status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all share this
#    string).
# 2) Narrow the "Status" / language-status columns:
#    - Overview sheet: columns E (zh-cn) and F (de-de)
#    - zh-cn sheet: column C (Status)
#    - de-de sheet: column C (Status)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update the status text ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status columns ---
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
